# Update "想去人数" (column F) figures on the "展览" (sheet1) and
# "全部类型" (sheet4) worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Map: row number -> new value for column F
$updates = @{
    2  = 1430
    3  = 7689
    6  = 35
    8  = 27
    9  = 5895
    10 = 152
    11 = 13
    12 = 28
    13 = 1803
    14 = 1312
    15 = 282
    16 = 399
    17 = 96
    18 = 5529
    19 = 65
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
